# Auto-generated Excel COM-interop edit script
# Applies cached numeric-value updates (market/profit data refresh)
# across the ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR leve-profit tables.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H34").Value = 3500
$ws.Range("I34").Value = 3500
$ws.Range("K34").Value = 3500
$ws.Range("M34").Value = -3297
$ws.Range("H36").Value = 3500
$ws.Range("I36").Value = 3500
$ws.Range("K36").Value = 3500
$ws.Range("M36").Value = -2785
$ws.Range("H112").Value = 6033.019
$ws.Range("I112").Value = 1900
$ws.Range("J112").Value = 6114.0586
$ws.Range("K112").Value = 5700
$ws.Range("L112").Value = 18342.1758
$ws.Range("M112").Value = -4592
$ws.Range("N112").Value = -20558.1758
$ws.Range("H129").Value = 1206.6842
$ws.Range("I129").Value = 591
$ws.Range("J129").Value = 1426.5714
$ws.Range("K129").Value = 1773
$ws.Range("L129").Value = 4279.7142
$ws.Range("M129").Value = 3227
$ws.Range("N129").Value = -14279.7142
$ws.Range("H132").Value = 5289.4814
$ws.Range("I132").Value = 4916.68
$ws.Range("J132").Value = 9949.5
$ws.Range("K132").Value = 14750.04
$ws.Range("L132").Value = 29848.5
$ws.Range("M132").Value = -12220.04
$ws.Range("N132").Value = -34908.5
$ws.Range("H133").Value = 0
$ws.Range("J133").Value = 0
$ws.Range("L133").Value = 0
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("H135").Value = 0
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("H136").Value = 30000
$ws.Range("J136").Value = 30000
$ws.Range("L136").Value = 30000
$ws.Range("N136").Value = -40200
$ws.Range("H137").Value = 4127.5
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 4127.5
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 12382.5
$ws.Range("N137").Value = -17482.5
$ws.Range("H138").Value = 1809.4445
$ws.Range("I138").Value = 2420.111
$ws.Range("J138").Value = 1565.1777
$ws.Range("K138").Value = 7260.333
$ws.Range("L138").Value = 4695.5331
$ws.Range("M138").Value = -2120.333
$ws.Range("N138").Value = -14975.5331
$ws.Range("N133").ClearContents()
$ws.Range("N134").ClearContents()
$ws.Range("M135").ClearContents()
$ws.Range("N135").ClearContents()
$ws.Range("M137").ClearContents()

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 528391.7
$ws.Range("I32").Value = 626253.3
$ws.Range("K32").Value = 626253.3
$ws.Range("M32").Value = -625966.3
$ws.Range("H61").Value = 4501.091
$ws.Range("I61").Value = 6506
$ws.Range("J61").Value = 4055.5557
$ws.Range("K61").Value = 6506
$ws.Range("L61").Value = 4055.5557
$ws.Range("M61").Value = -6294
$ws.Range("N61").Value = -4479.5557
$ws.Range("H97").Value = 1060.8889
$ws.Range("I97").Value = 943.5
$ws.Range("J97").Value = 2000
$ws.Range("K97").Value = 943.5
$ws.Range("L97").Value = 2000
$ws.Range("M97").Value = -447.5
$ws.Range("N97").Value = -2992
$ws.Range("H102").Value = 1935.1765
$ws.Range("I102").Value = 1915.3077
$ws.Range("K102").Value = 1915.3077
$ws.Range("M102").Value = -293.3077000000001
$ws.Range("H132").Value = 7813.75
$ws.Range("I132").Value = 11128
$ws.Range("J132").Value = 4499.5
$ws.Range("K132").Value = 33384
$ws.Range("L132").Value = 13498.5
$ws.Range("M132").Value = -30854
$ws.Range("N132").Value = -18558.5
$ws.Range("H136").Value = 4501.091
$ws.Range("I136").Value = 6506
$ws.Range("J136").Value = 4055.5557
$ws.Range("K136").Value = 19518
$ws.Range("L136").Value = 12166.6671
$ws.Range("M136").Value = -16968
$ws.Range("N136").Value = -17266.6671

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 58826184
$ws.Range("I86").Value = 76925860
$ws.Range("J86").Value = 2250
$ws.Range("K86").Value = 76925860
$ws.Range("L86").Value = 2250
$ws.Range("M86").Value = -76924737
$ws.Range("N86").Value = -4496
$ws.Range("H89").Value = 58826184
$ws.Range("I89").Value = 76925860
$ws.Range("J89").Value = 2250
$ws.Range("K89").Value = 384629300
$ws.Range("L89").Value = 11250
$ws.Range("M89").Value = -384623684
$ws.Range("N89").Value = -22482
$ws.Range("H94").Value = 1071
$ws.Range("I94").Value = 696.3333
$ws.Range("J94").Value = 1295.8
$ws.Range("K94").Value = 696.3333
$ws.Range("L94").Value = 1295.8
$ws.Range("M94").Value = -245.3333
$ws.Range("N94").Value = -2197.8
$ws.Range("H105").Value = 8930686
$ws.Range("I105").Value = 10418683
$ws.Range("J105").Value = 2700
$ws.Range("K105").Value = 10418683
$ws.Range("L105").Value = 2700
$ws.Range("M105").Value = -10416936
$ws.Range("N105").Value = -6194
$ws.Range("H134").Value = 3800
$ws.Range("I134").Value = 3866.6667
$ws.Range("J134").Value = 3666.6667
$ws.Range("K134").Value = 11600.0001
$ws.Range("L134").Value = 11000.0001
$ws.Range("M134").Value = -9065.000100000001
$ws.Range("N134").Value = -16070.0001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H105").Value = 966.6667
$ws.Range("I105").Value = 1000
$ws.Range("J105").Value = 900
$ws.Range("K105").Value = 1000
$ws.Range("L105").Value = 900
$ws.Range("M105").Value = 747
$ws.Range("N105").Value = -4394
$ws.Range("H112").Value = 35333.332
$ws.Range("J112").Value = 35333.332
$ws.Range("L112").Value = 35333.332
$ws.Range("N112").Value = -38287.332
$ws.Range("H119").Value = 0
$ws.Range("J119").Value = 0
$ws.Range("L119").Value = 0
$ws.Range("H132").Value = 8774109
$ws.Range("I132").Value = 1541
$ws.Range("J132").Value = 23812796
$ws.Range("K132").Value = 4623
$ws.Range("L132").Value = 71438388
$ws.Range("M132").Value = -2093
$ws.Range("N132").Value = -71443448
$ws.Range("H134").Value = 3752
$ws.Range("I134").Value = 4537.3335
$ws.Range("J134").Value = 2966.6667
$ws.Range("K134").Value = 13612.0005
$ws.Range("L134").Value = 8900.000100000001
$ws.Range("M134").Value = -11077.0005
$ws.Range("N134").Value = -13970.0001
$ws.Range("N119").ClearContents()

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H3").Value = 19956.5
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = 19956.5
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = 59869.5
$ws.Range("N3").Value = -60093.5
$ws.Range("H5").Value = 976.8889
$ws.Range("I5").Value = 852.875
$ws.Range("J5").Value = 1969
$ws.Range("K5").Value = 2558.625
$ws.Range("L5").Value = 5907
$ws.Range("M5").Value = -2446.625
$ws.Range("N5").Value = -6131
$ws.Range("H60").Value = 2898.2456
$ws.Range("I60").Value = 400
$ws.Range("J60").Value = 3037.037
$ws.Range("K60").Value = 1200
$ws.Range("L60").Value = 9111.110999999999
$ws.Range("M60").Value = -949
$ws.Range("N60").Value = -9613.110999999999
$ws.Range("H113").Value = 823.34375
$ws.Range("I113").Value = 452.16666
$ws.Range("J113").Value = 1300.5714
$ws.Range("K113").Value = 1356.49998
$ws.Range("L113").Value = 3901.7142
$ws.Range("M113").Value = 813.5000199999999
$ws.Range("N113").Value = -8241.7142
$ws.Range("H131").Value = 1024.3448
$ws.Range("J131").Value = 1139.08
$ws.Range("L131").Value = 3417.24
$ws.Range("N131").Value = -13497.24
$ws.Range("H135").Value = 976.8889
$ws.Range("I135").Value = 852.875
$ws.Range("J135").Value = 1969
$ws.Range("K135").Value = 7675.875
$ws.Range("L135").Value = 17721
$ws.Range("M135").Value = -5140.875
$ws.Range("N135").Value = -22791
$ws.Range("M3").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3994.0232
$ws.Range("I122").Value = 1388.4375
$ws.Range("J122").Value = 5538.074
$ws.Range("K122").Value = 4165.3125
$ws.Range("L122").Value = 16614.222
$ws.Range("M122").Value = -1715.3125
$ws.Range("N122").Value = -21514.222
$ws.Range("H132").Value = 4187.375
$ws.Range("I132").Value = 3916.6667
$ws.Range("J132").Value = 4999.5
$ws.Range("K132").Value = 11750.0001
$ws.Range("L132").Value = 14998.5
$ws.Range("M132").Value = -9220.000100000001
$ws.Range("N132").Value = -20058.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 29375.5
$ws.Range("J2").Value = 29375.5
$ws.Range("L2").Value = 29375.5
$ws.Range("N2").Value = -29599.5
$ws.Range("H82").Value = 3250.5
$ws.Range("I82").Value = 3250.5
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 3250.5
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -2889.5
$ws.Range("H85").Value = 3250.5
$ws.Range("I85").Value = 3250.5
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 3250.5
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -2002.5
$ws.Range("H93").Value = 8077.1113
$ws.Range("I93").Value = 10390.417
$ws.Range("J93").Value = 3450.5
$ws.Range("K93").Value = 10390.417
$ws.Range("L93").Value = 3450.5
$ws.Range("M93").Value = -9142.416999999999
$ws.Range("N93").Value = -5946.5
$ws.Range("H140").Value = 82214.5
$ws.Range("J140").Value = 82214.5
$ws.Range("L140").Value = 82214.5
$ws.Range("N140").Value = -92574.5
$ws.Range("N82").ClearContents()
$ws.Range("N85").ClearContents()

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 52639.668
$ws.Range("J46").Value = 52639.668
$ws.Range("L46").Value = 52639.668
$ws.Range("N46").Value = -53101.668
$ws.Range("H100").Value = 1139.8
$ws.Range("I100").Value = 699.6667
$ws.Range("J100").Value = 1800
$ws.Range("K100").Value = 1399.3334
$ws.Range("L100").Value = 3600
$ws.Range("M100").Value = -858.3334
$ws.Range("N100").Value = -4682
$ws.Range("H123").Value = 31742.834
$ws.Range("J123").Value = 53095.668
$ws.Range("L123").Value = 53095.668
$ws.Range("N123").Value = -62895.668
$ws.Range("H125").Value = 49400
$ws.Range("J125").Value = 49400
$ws.Range("L125").Value = 49400
$ws.Range("N125").Value = -59240
$ws.Range("H134").Value = 52639.668
$ws.Range("J134").Value = 52639.668
$ws.Range("L134").Value = 157919.004
$ws.Range("N134").Value = -162989.004
